$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.439.58'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.79%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.793.82'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.34%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.06%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.558'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.92%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.85'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.12%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.297'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.68%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0694'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.94%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0950'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.45%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.052.61'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.26%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.797.32'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.57%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.07'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.62%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.637'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.62%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '34.439.77'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.88%  '

# Row 17
$ws.Range('E17').Value = '  +2.36%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.89'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.93%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '247.16'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.44%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.84%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.22'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.90%  '

# Row 22
$ws.Range('E22').Value = '  -0.04%  '

# Row 23
$ws.Range('E23').Value = '  +1.66%  '

# Row 24
$ws.Range('E24').Value = '  +0.97%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '164.62'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.23%  '

# Row 26
$ws.Range('E26').Value = '  +1.17%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.52'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.21%  '

# Row 28
$ws.Range('E28').Value = '  +2.66%  '

# Row 29
$ws.Range('E29').Value = '  -0.14%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.82'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.20%  '

# Row 31
$ws.Range('E31').Value = '  +0.00%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0522'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.80%  '

# Row 33
$ws.Range('E33').Value = '  +7.07%  '

# Row 34
$ws.Range('E34').Value = '  +1.05%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.431.28'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.85%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.59'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.31%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.671'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.52%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.80%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0192'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.14%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '84.77'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.62%  '

# Row 41
$ws.Range('E41').Value = '  +1.03%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.940'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.27%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.74'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.16%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.56'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.69%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0524'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.93%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.09'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.78%  '

# Row 47
$ws.Range('E47').Value = '  +0.18%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.949.04'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.02%  '

# Row 49
$ws.Range('E49').Value = '  -0.05%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0131'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.16%  '

# Row 51
$ws.Range('E51').Value = '  -0.04%  '

